$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.292.06"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.171.80"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.02"
$ws.Range("E5").Value = "  +5.76%  "

$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.27"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.41"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("D14").Value = "2.497.48"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.09"
$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").Value = "2.152.95"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("E17").Value = "  -3.04%  "

$ws.Range("D18").Value = "42.186.08"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("E19").Value = "  -1.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.38"
$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  +0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.14"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.13"
$ws.Range("E23").Value = "  +3.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.39"
$ws.Range("E24").Value = "  -7.15%  "

$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("E26").Value = "  -3.89%  "

$ws.Range("E27").Value = "  +1.76%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  +6.44%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.62"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.51"
$ws.Range("E31").Value = "  +10.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.94"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.08"
$ws.Range("E34").Value = "  -4.31%  "

$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("E37").Value = "  -4.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0332"
$ws.Range("E38").Value = "  +5.88%  "

$ws.Range("E39").Value = "  -2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.66"
$ws.Range("E40").Value = "  -5.63%  "

$ws.Range("E41").Value = "  +1.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "58.99"
$ws.Range("E42").Value = "  -1.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.12"
$ws.Range("E43").Value = "  -5.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.55"
$ws.Range("E44").Value = "  +4.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.467"
$ws.Range("E45").Value = "  +12.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0969"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("E47").Value = "  -3.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("E48").Value = "  +7.63%  "

$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  +0.39%  "

Write-Output "Update complete"